$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7174.4287
$ws.Range("J40").Value = 8714
$ws.Range("L40").Value = 8714
$ws.Range("N40").Value = -9064
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 10039.833
$ws.Range("I132").Value = 10134.363
$ws.Range("K132").Value = 30403.089
$ws.Range("M132").Value = -27873.089
$ws.Range("H138").Value = 1536.8
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2816.0715
$ws.Range("I45").Value = 1475
$ws.Range("K45").Value = 1475
$ws.Range("M45").Value = -1098
$ws.Range("H69").Value = 249999.5
$ws.Range("J69").Value = 249999.5
$ws.Range("L69").Value = 249999.5
$ws.Range("N69").Value = -251497.5
$ws.Range("H72").Value = 249999.5
$ws.Range("J72").Value = 249999.5
$ws.Range("L72").Value = 749998.5
$ws.Range("N72").Value = -757486.5
$ws.Range("H119").Value = 100001
$ws.Range("J119").Value = 100001
$ws.Range("L119").Value = 100001
$ws.Range("N119").Value = -109677

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 12636
$ws.Range("J54").Value = 21000
$ws.Range("L54").Value = 21000
$ws.Range("N54").Value = -21968
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 34844.375
$ws.Range("J41").Value = 35057.5
$ws.Range("L41").Value = 35057.5
$ws.Range("N41").Value = -35913.5
$ws.Range("H50").Value = 97000
$ws.Range("J50").Value = 97000
$ws.Range("L50").Value = 97000
$ws.Range("N50").Value = -98250
$ws.Range("H56").Value = 500
$ws.Range("I56").Value = 500
$ws.Range("K56").Value = 500
$ws.Range("M56").Value = 345
$ws.Range("H59").Value = 329280000
$ws.Range("J59").Value = 493876540
$ws.Range("L59").Value = 493876540
$ws.Range("N59").Value = -493878830
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H86").Value = 3016.4
$ws.Range("I86").Value = 2860.6667
$ws.Range("J86").Value = 3250
$ws.Range("K86").Value = 2860.6667
$ws.Range("L86").Value = 3250
$ws.Range("M86").Value = -1737.6667
$ws.Range("N86").Value = -5496
$ws.Range("H89").Value = 3016.4
$ws.Range("I89").Value = 2860.6667
$ws.Range("J89").Value = 3250
$ws.Range("K89").Value = 14303.3335
$ws.Range("L89").Value = 16250
$ws.Range("M89").Value = -8687.333500000001
$ws.Range("N89").Value = -27482
$ws.Range("H99").Value = 15000
$ws.Range("I99").Value = 15000
$ws.Range("K99").Value = 15000
$ws.Range("M99").Value = -13502
$ws.Range("H103").Value = 18233.75
$ws.Range("I103").Value = 18233.75
$ws.Range("K103").Value = 18233.75
$ws.Range("M103").Value = -17061.75
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("H126").Value = 15000
$ws.Range("I126").Value = 15000
$ws.Range("K126").Value = 45000
$ws.Range("M126").Value = -42530
$ws.Range("H134").Value = 4000
$ws.Range("I134").Value = 4000
$ws.Range("K134").Value = 12000
$ws.Range("M134").Value = -9465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12799.5
$ws.Range("I3").Value = 12799
$ws.Range("K3").Value = 38397
$ws.Range("M3").Value = -38285
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 15000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -16372
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 45000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -51864
$ws.Range("H131").Value = 1943.6316
$ws.Range("I131").Value = 629.8333
$ws.Range("K131").Value = 1889.4999
$ws.Range("M131").Value = 3150.5001
$ws.Range("H138").Value = 6006.25
$ws.Range("J138").Value = 8710
$ws.Range("L138").Value = 26130
$ws.Range("N138").Value = -36410

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 65007
$ws.Range("J24").Value = 65007
$ws.Range("L24").Value = 65007
$ws.Range("N24").Value = -65353
$ws.Range("H68").Value = 45268
$ws.Range("I68").Value = 45268
$ws.Range("K68").Value = 45268
$ws.Range("M68").Value = -44457
$ws.Range("H71").Value = 45268
$ws.Range("I71").Value = 45268
$ws.Range("K71").Value = 135804
$ws.Range("M71").Value = -131748

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6960.7144
$ws.Range("J46").Value = 8875
$ws.Range("L46").Value = 8875
$ws.Range("N46").Value = -9251
$ws.Range("H61").Value = 6325.2856
$ws.Range("I61").Value = 2500
$ws.Range("J61").Value = 7855.4
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 7855.4
$ws.Range("M61").Value = -2298
$ws.Range("N61").Value = -8259.4
$ws.Range("H113").Value = 6325.2856
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 7855.4
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 7855.4
$ws.Range("M113").Value = -330
$ws.Range("N113").Value = -12195.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H96").Value = 1399.5
$ws.Range("I96").Value = 1399.5
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1399.5
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -26.5
$ws.Range("N96").ClearContents()
$ws.Range("H113").Value = 724.5789
$ws.Range("I113").Value = 560.4167
$ws.Range("J113").Value = 1006
$ws.Range("K113").Value = 1681.2501
$ws.Range("L113").Value = 3018
$ws.Range("M113").Value = 488.7499
$ws.Range("N113").Value = -7358
